$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply formatting first (copy from existing rows with matching style), so it
# does not interfere with the shared-string allocation order below. ---
$ws.Range("A46:E46").Copy()
$ws.Range("A56:E56").PasteSpecial(-4122)
$ws.Range("A57:E57").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A58:E58").PasteSpecial(-4122)
$ws.Range("A60:E60").PasteSpecial(-4122)
$ws.Range("A4:E4").Copy()
$ws.Range("A59:E59").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row heights (auto-computed by real Excel when it wraps this text at the
# sheet's column widths; pinned explicitly here to match).
$ws.Rows.Item(56).RowHeight = 43.2
$ws.Rows.Item(57).RowHeight = 43.2
$ws.Rows.Item(58).RowHeight = 43.2
$ws.Rows.Item(59).RowHeight = 21.6
$ws.Rows.Item(60).RowHeight = 43.2

# --- Set cell values. Order matches the source workbook's shared-string table
# (EN column(s) of a block, then the SCRIPT file ref, then RU column(s), then the
# "converted" column(s)) so newly appended <si> entries land on the same indices. ---
$ws.Cells.Item(56,3).Value = " I got to see the one and only\nTeam [CS:X]Charm[CR] with my own eyes...[K] That\'s more\nthan enough for me!"
$ws.Cells.Item(56,1).Value = "SCRIPT/G01P03A/us0401.ssb"
$ws.Cells.Item(56,4).Value = " Я собственными глазами видела\nту самую Команду [CS:X]Шарм[CR]...[K] Мне этого более\nчем достаточно!"
$ws.Cells.Item(56,5).Value = " Ÿ òïáòóâåîîúíé ãìàèàíé âéäåìà\nóô òàíôý Ëïíàîäô [CS:X]Šàñí[CR]...[K] Íîå üóïãï áïìåå\nœåí äïòóàóïœîï!"
$ws.Cells.Item(57,3).Value = " I hope you manage to bring\n[CS:N]Drowzee[CR] back!"
$ws.Cells.Item(57,1).Value = "SCRIPT/G01P03A/us2006.ssb"
$ws.Cells.Item(57,4).Value = " Надеюсь, вы сумеете привести\nсюда [CS:N]Дроузи[CR]!"
$ws.Cells.Item(57,5).Value = " Îàäåýòû, âú òôíååóå ðñéâåòóé\nòýäà [CS:N]Äñïôèé[CR]!"
$ws.Cells.Item(58,3).Value = " Oh?[K] You can go into [CS:N]Azurill[CR]\'s\ndream?[K] That\'s wonderful news!"
$ws.Cells.Item(59,3).Value = " I hope you\'ll get to the root of\nthe problem!"
$ws.Cells.Item(58,1).Value = "SCRIPT/G01P03A/us2009.ssb"
$ws.Cells.Item(58,4).Value = " Что?[K] Вы сможете попасть в сон\n[CS:N]Азурилла[CR]?[K] Это же прекрасно!"
$ws.Cells.Item(59,4).Value = " Надеюсь, вы доберётесь до истины!"
$ws.Cells.Item(58,5).Value = " Œóï?[K] Âú òíïçåóå ðïðàòóû â òïî\n[CS:N]Àèôñéììà[CR]?[K] Üóï çå ðñåëñàòîï!"
$ws.Cells.Item(59,5).Value = " Îàäåýòû, âú äïáåñæóåòû äï éòóéîú!"
$ws.Cells.Item(60,3).Value = " Good luck!"
$ws.Cells.Item(60,1).Value = "SCRIPT/G01P03A/us2201.ssb"
$ws.Cells.Item(60,4).Value = " Удачи вам!"
$ws.Cells.Item(60,5).Value = " Ôäàœé âàí!"

# Row-number column (plain integers, not shared strings).
$ws.Cells.Item(56,2).Value = 206
$ws.Cells.Item(57,2).Value = 151
$ws.Cells.Item(58,2).Value = 123
$ws.Cells.Item(59,2).Value = 126
$ws.Cells.Item(60,2).Value = 68

# Match the author's final selection/view state.
$ws.Range("D58").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 56
$win.ScrollColumn = 1
